$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows for Nonantola, aggiornamento a 9/09 compreso (through 2021-09-09).
$data = @(
    @(367, 44441, 6, 21, 131.4965560425798),
    @(368, 44442, 1, 19, 118.9730745147151),
    @(369, 44443, 0, 17, 106.4495929868503),
    @(370, 44444, 0, 14, 87.66437069505322),
    @(371, 44445, 1, 9, 56.35566687539136),
    @(372, 44446, 0, 9, 56.35566687539136),
    @(373, 44447, 0, 8, 50.09392611145898),
    @(374, 44448, 0, 2, 12.52348152786475)
)

# Column A carries the same date-number style as the rest of the column
# (bordered, bold, centered, custom date format) -- copy it down first.
$ws.Range("A366").Copy()
$ws.Range("A367:A374").PasteSpecial(-4122)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}
